# Commit: "add few fields and error report"
# Updates the product rows (2-4) with corrected catalog data and adds a
# new 5th product row (previously just blank date-formatted cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = "5116877-32863368197"
$ws.Range("B2").Value = "NYLON RUBBER BAND Pkt 500gm"
$ws.Range("C2").Value = "rubber bands"
$ws.Range("D2").Value = "Ekam"
$ws.Range("E2").Value = "http://mkp.gem.gov.in/rubber-bands/nylon-rubber-band-pkt-500gm/p-5116877-32863368197-cat.html"
$ws.Range("F2").Value = "https://admin-mkp.gem.gov.in/#!/catalog/new?bnid=home_offi_of45811733_fast_rubb&gem_catalog_id=5116877-32863368197"
$ws.Range("G2").Value = "ST-GI12/A472-20-25"
$ws.Range("H2").Value = "Genaric"
$ws.Range("I2").Value = 44166
$ws.Range("J2").Value = 44166
$ws.Range("K2").Value = 46022
$ws.Range("L2").Value = "India"
$ws.Range("N2").Value = 4016
$ws.Range("O2").Value = 312
$ws.Range("P2").Value = 279.99
$ws.Range("S2").Value = "all"
$ws.Range("T2").Value = 51415
$ws.Range("U2").Value = 6
$ws.Range("V2").Value = 3

# --- Row 3 ---
$ws.Range("A3").Value = "5116877-15821908934"
$ws.Range("B3").Value = "GOOD MAKE RUBBER BAND"
$ws.Range("C3").Value = "rubber bands"
$ws.Range("D3").Value = "Good Make"
$ws.Range("E3").Value = "http://mkp.gem.gov.in/rubber-bands/rubber-band-big-size/p-5116877-15821908934-cat.html"
$ws.Range("F3").Value = "https://admin-mkp.gem.gov.in/#!/catalog/new?bnid=home_offi_of45811733_fast_rubb&gem_catalog_id=5116877-15821908934"
$ws.Range("G3").Value = "ST-GI12/A472-20-25"
$ws.Range("H3").Value = "Genaric"
$ws.Range("I3").Value = 44166
$ws.Range("J3").Value = 44166
$ws.Range("K3").Value = 46022
$ws.Range("L3").Value = "India"
$ws.Range("N3").Value = 4016
$ws.Range("O3").Value = 480
$ws.Range("P3").Value = 398.99
$ws.Range("S3").Value = "all"
$ws.Range("T3").Value = 51417
$ws.Range("U3").Value = 4
$ws.Range("V3").Value = 3

# --- Row 4 ---
$ws.Range("A4").Value = "5116877-17154744803"
$ws.Range("B4").Value = "R-73947"
$ws.Range("C4").Value = "rubber bands"
$ws.Range("D4").Value = "Good Make"
$ws.Range("E4").Value = "http://mkp.gem.gov.in/office-equipment-accessories-supplies/rubber-bands/p-5116877-17154744803-cat.html"
$ws.Range("F4").Value = "https://admin-mkp.gem.gov.in/#!/catalog/new?bnid=home_offi_of45811733_fast_rubb&gem_catalog_id=5116877-17154744803"
$ws.Range("G4").Value = "ST-GI12/A472-20-25"
$ws.Range("H4").Value = "Genaric"
$ws.Range("I4").Value = 44166
$ws.Range("J4").Value = 44166
$ws.Range("K4").Value = 46022
$ws.Range("L4").Value = "India"
$ws.Range("N4").Value = 4016
$ws.Range("O4").Value = 430
$ws.Range("P4").Value = 380.99
$ws.Range("S4").Value = "all"
$ws.Range("T4").Value = 51418
$ws.Range("U4").Value = 4
$ws.Range("V4").Value = 3

# --- Row 5 (new product row; previously only had blank date cells I5:K5) ---
$ws.Range("A5").Value = "5116877-92981260387"
$ws.Range("B5").Value = "RUBBER BAND GOOD MAKE"
$ws.Range("C5").Value = "rubber bands"
$ws.Range("D5").Value = "Good Make"
$ws.Range("E5").Value = "http://mkp.gem.gov.in/rubber-bands/rubber-band-good-make/p-5116877-92981260387-cat.html"
$ws.Range("F5").Value = "https://admin-mkp.gem.gov.in/#!/catalog/new?bnid=home_offi_of45811733_fast_rubb&gem_catalog_id=5116877-92981260387"
$ws.Range("G5").Value = "ST-GI12/A472-20-25"
$ws.Range("H5").Value = "Genaric"
$ws.Range("I5").Value = 44166
$ws.Range("J5").Value = 44166
$ws.Range("K5").Value = 46022
$ws.Range("L5").Value = "India"
$ws.Range("N5").Value = 4016
$ws.Range("O5").Value = 430
$ws.Range("P5").Value = 380.99
$ws.Range("S5").Value = "all"
$ws.Range("T5").Value = 51420
$ws.Range("U5").Value = 4
$ws.Range("V5").Value = 3

# --- Selection moved to B9 ---
$ws.Range("B9").Select()
